$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K7").Value  = 0.2970525035592049
$ws.Range("J8").Value  = 0.2858677898194339
$ws.Range("I9").Value  = 0.2775335613519331
$ws.Range("H10").Value = 0.2743085116504074
$ws.Range("G11").Value = 0.2534447081011285
$ws.Range("F12").Value = 0.2766837437271186
$ws.Range("E13").Value = 0.2867219094086165
$ws.Range("D14").Value = 0.1751453671933744
$ws.Range("C15").Value = 0.1965658720679752
$ws.Range("B16").Value = 0.4328090033804217
